# "also doing nr topics 20 reg fitting"
# Add a new results row (row 6) for another model run: sample_model=2500,
# sample_companies=300, type=reg, nr topics=20.
# Also jot a note about an intermediate result in S4, and leave the
# active selection on G8 (where the user was about to continue working).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note added next to the row-4 "TODOs" note.
$ws.Range("S4").Value = ".-> 6700 nachdem etwas toleranz für die industry calls eingebaut wurde!S9"

# New run results in row 6.
$ws.Range("B6").Value = 2500
$ws.Range("C6").Value = 300
$ws.Range("D6").Value = "reg"
$ws.Range("E6").Value = 20

# Leave selection where the user left off.
$ws.Range("G8").Select()
